$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-02-23 Friday"

# Update each table cell value (row-major order, 20 rows x 5 columns)
$t = $d.Tables.Item(1)
$newValues = @(
    "16+77=93",
    "91-53=38",
    "33-28=5",
    "17+37=54",
    "90-57=33",
    "40-28=12",
    "31-16=15",
    "55+27=82",
    "61-28=33",
    "16+28=44",
    "57+14=71",
    "48+3=51",
    "84-75=9",
    "27+18=45",
    "90-24=66",
    "19+19=38",
    "52-18=34",
    "38+28=66",
    "39+5=44",
    "8+14=22",
    "19+55=74",
    "28+34=62",
    "74-36=38",
    "37+19=56",
    "70-22=48",
    "34+48=82",
    "53-5=48",
    "8+4=12",
    "47+38=85",
    "5+79=84",
    "63-36=27",
    "73-56=17",
    "14+19=33",
    "43-28=15",
    "7+29=36",
    "43-14=29",
    "28+66=94",
    "61-53=8",
    "40-13=27",
    "90-15=75",
    "18+15=33",
    "27+17=44",
    "28+65=93",
    "60-2=58",
    "81-46=35",
    "26+56=82",
    "61-15=46",
    "64+28=92",
    "65-17=48",
    "6+45=51",
    "9+82=91",
    "28+33=61",
    "6+35=41",
    "46+7=53",
    "31-15=16",
    "52-3=49",
    "45+16=61",
    "9+25=34",
    "65-7=58",
    "41-22=19",
    "37+56=93",
    "77+4=81",
    "18+48=66",
    "29+7=36",
    "18+73=91",
    "59+16=75",
    "8+65=73",
    "34-5=29",
    "45+8=53",
    "61-36=25",
    "77+19=96",
    "74-29=45",
    "7+68=75",
    "9+27=36",
    "72-15=57",
    "60-1=59",
    "91-2=89",
    "14+8=22",
    "54-8=46",
    "94-36=58",
    "49+9=58",
    "37+45=82",
    "40-39=1",
    "50-33=17",
    "7+79=86",
    "37+57=94",
    "40-2=38",
    "4+59=63",
    "35+18=53",
    "39+34=73",
    "28+66=94",
    "46+26=72",
    "43-17=26",
    "19+22=41",
    "56+6=62",
    "45+28=73",
    "52-33=19",
    "80-56=24",
    "8+43=51",
    "90-51=39"
)

$idx = 0
for ($row = 1; $row -le 20; $row++) {
    for ($col = 1; $col -le 5; $col++) {
        $t.Cell($row, $col).Range.Text = $newValues[$idx]
        $idx++
    }
}
